$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = "Cole Garrett"
$ws.Range("C2").Value = "'1222222250"
$ws.Range("C2").Style = "Normal"
$ws.Range("F2").Value = "Pediatric Emergency Medicine"
$ws.Range("H2").Value = "'08/01/2025"
$ws.Range("H2").Style = "Normal"
$ws.Range("I2").Value = "'08/01/2025"
$ws.Range("I2").Style = "Normal"
$ws.Range("K2").Value = "'"
$ws.Range("K2").Style = "Normal"

# Row 3 updates
$ws.Range("D3").Value = "'458888885"
$ws.Range("D3").Style = "Normal"
$ws.Range("J3").Value = "PPG#'s, Medicare, Commercial HMO"
$ws.Range("K3").Value = "'"
$ws.Range("K3").Style = "Normal"
$ws.Range("L3").Value = "Mercian Medical Group - 1014"

# Row 4 updates
$ws.Range("S4").Value = "prajay.sapkale@hilabs.com"
